$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D5").Value = "2016-01-28 09:43:56"
$wsZh.Range("G5").Value = "2016-01-28 09:44:38"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D5").Value = "2016-01-28 09:44:10"
$wsDe.Range("G5").Value = "2016-01-28 09:45:01"
